# Refatorando o consolidador para modelo ETL
# Update absenteeism data rows 2-11 with new values produced by the ETL pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date serial epoch (1900 date system, with the historical 1900 leap-year bug)
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

# Row 2
$ws.Cells.Item(2, 1).Value = 96368
$ws.Cells.Item(2, 2).Value = "André Azevedo"
$ws.Cells.Item(2, 3).Value = "TI"
$ws.Cells.Item(2, 4).Value = "Doença"
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = $epoch.AddDays(45091)
$ws.Cells.Item(2, 7).Value = 9383.83

# Row 3
$ws.Cells.Item(3, 1).Value = 18745
$ws.Cells.Item(3, 2).Value = "Diogo Monteiro"
$ws.Cells.Item(3, 3).Value = "Recursos Humanos"
$ws.Cells.Item(3, 4).Value = "Viagem de negócios"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = $epoch.AddDays(45093)
$ws.Cells.Item(3, 7).Value = 5661.68

# Row 4
$ws.Cells.Item(4, 1).Value = 63
$ws.Cells.Item(4, 2).Value = "Dr. Benício Monteiro"
$ws.Cells.Item(4, 3).Value = "Vendas"
$ws.Cells.Item(4, 4).Value = "Consulta médica"
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = $epoch.AddDays(45094)
$ws.Cells.Item(4, 7).Value = 8202.51

# Row 5
$ws.Cells.Item(5, 1).Value = 32942
$ws.Cells.Item(5, 2).Value = "Luiz Fernando Alves"
$ws.Cells.Item(5, 3).Value = "Vendas"
$ws.Cells.Item(5, 4).Value = "Viagem de negócios"
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = $epoch.AddDays(45091)
$ws.Cells.Item(5, 7).Value = 3209.15

# Row 6
$ws.Cells.Item(6, 1).Value = 34203
$ws.Cells.Item(6, 2).Value = "Kevin Araújo"
$ws.Cells.Item(6, 3).Value = "Marketing"
$ws.Cells.Item(6, 4).Value = "Viagem de negócios"
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = $epoch.AddDays(45105)
$ws.Cells.Item(6, 7).Value = 7483.92

# Row 7
$ws.Cells.Item(7, 1).Value = 23621
$ws.Cells.Item(7, 2).Value = "Vitor Barbosa"
$ws.Cells.Item(7, 3).Value = "Marketing"
$ws.Cells.Item(7, 4).Value = "Doença"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = $epoch.AddDays(45088)
$ws.Cells.Item(7, 7).Value = 8212.55

# Row 8
$ws.Cells.Item(8, 1).Value = 42736
$ws.Cells.Item(8, 2).Value = "Marcela Nascimento"
$ws.Cells.Item(8, 3).Value = "P&D"
$ws.Cells.Item(8, 4).Value = "Consulta médica"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = $epoch.AddDays(45095)
$ws.Cells.Item(8, 7).Value = 2612.76

# Row 9
$ws.Cells.Item(9, 1).Value = 15541
$ws.Cells.Item(9, 2).Value = "Maria Sophia Cunha"
$ws.Cells.Item(9, 3).Value = "TI"
$ws.Cells.Item(9, 4).Value = "Outros"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = $epoch.AddDays(45098)
$ws.Cells.Item(9, 7).Value = 3946.97

# Row 10
$ws.Cells.Item(10, 1).Value = 96293
$ws.Cells.Item(10, 2).Value = "Dr. Diogo Barbosa"
$ws.Cells.Item(10, 3).Value = "Marketing"
$ws.Cells.Item(10, 4).Value = "Doença"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = $epoch.AddDays(45083)
$ws.Cells.Item(10, 7).Value = 3201.48

# Row 11
$ws.Cells.Item(11, 1).Value = 90945
$ws.Cells.Item(11, 2).Value = "Sophia da Cruz"
$ws.Cells.Item(11, 3).Value = "Engenharia"
$ws.Cells.Item(11, 4).Value = "Outros"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = $epoch.AddDays(45099)
$ws.Cells.Item(11, 7).Value = 2719.45
